# Regenerate the "K" (strikeouts) column (G) values for wittgren_nick.xlsx
# save_data sheet. The source data was regenerated to use actual strikeout
# counts (K) in place of the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 1
    31 = 1
    32 = 2
    33 = 1
    34 = 1
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
